# paises.xlsx -- refresh COVID-19 country snapshot (27 Mar 2020, 07:12 -> 07:42)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell: "Datos actualizados a ..." timestamp
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 07:42"

# Row 4: Estados Unidos -- updated totals
$ws.Cells.Item(4,2).Value = 85604
$ws.Cells.Item(4,3).Value = 169
$ws.Cells.Item(4,5).Value = 82435
$ws.Cells.Item(4,7).Value = 6
$ws.Cells.Item(4,8).Value = 1301

# Row 21: Australia -- updated totals
$ws.Cells.Item(21,2).Value = 3143
$ws.Cells.Item(21,3).Value = 93
$ws.Cells.Item(21,5).Value = 2960

# Row 33: Pakistan -- row 33/34 swap with Polonia -- Pakistan overtakes it, new Pakistan stats
$ws.Cells.Item(33,1).Value = "Pakistan"
$ws.Cells.Item(33,2).Value = 1235
$ws.Cells.Item(33,3).Value = 34
$ws.Cells.Item(33,4).Value = 23
$ws.Cells.Item(33,5).Value = 1203
$ws.Cells.Item(33,6).Value = 7
$ws.Cells.Item(33,8).Value = 9

# Row 34: Polonia -- row 33/34 swap with Pakistan -- Polonia drops a rank, stats unchanged
$ws.Cells.Item(34,1).Value = "Polonia"
$ws.Cells.Item(34,2).Value = 1221
$ws.Cells.Item(34,3).Value = 0
$ws.Cells.Item(34,4).Value = 7
$ws.Cells.Item(34,5).Value = 1198
$ws.Cells.Item(34,6).Value = 3
$ws.Cells.Item(34,8).Value = 16

# Row 39: Sudafrica -- updated deaths-today
$ws.Cells.Item(39,6).Value = 7

# Row 44: India -- updated totals
$ws.Cells.Item(44,2).Value = 747
$ws.Cells.Item(44,3).Value = 20
$ws.Cells.Item(44,5).Value = 661

# Row 66: Lituania -- rows 66-72 reshuffle -- Lituania jumps to row 66 with new stats
$ws.Cells.Item(66,1).Value = "Lituania"
$ws.Cells.Item(66,2).Value = 344
$ws.Cells.Item(66,3).Value = 45
$ws.Cells.Item(66,4).Value = 1
$ws.Cells.Item(66,5).Value = 339
$ws.Cells.Item(66,6).Value = 1
$ws.Cells.Item(66,8).Value = 4

# Row 67: Emiratos Arabes Unidos -- rows 66-72 reshuffle -- shifts to row 67, stats unchanged
$ws.Cells.Item(67,1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(67,2).Value = 333
$ws.Cells.Item(67,3).Value = 0
$ws.Cells.Item(67,4).Value = 52
$ws.Cells.Item(67,5).Value = 279
$ws.Cells.Item(67,6).Value = 2
$ws.Cells.Item(67,8).Value = 2

# Row 68: Armenia -- rows 66-72 reshuffle -- Armenia jumps to row 68 with new stats
$ws.Cells.Item(68,1).Value = "Armenia"
$ws.Cells.Item(68,2).Value = 329
$ws.Cells.Item(68,3).Value = 39
$ws.Cells.Item(68,4).Value = 18
$ws.Cells.Item(68,5).Value = 310
$ws.Cells.Item(68,6).Value = 6
$ws.Cells.Item(68,8).Value = 1

# Row 69: Hungria -- rows 66-72 reshuffle -- shifts to row 69, stats unchanged
$ws.Cells.Item(69,1).Value = "Hungria"
$ws.Cells.Item(69,2).Value = 300
$ws.Cells.Item(69,3).Value = 39
$ws.Cells.Item(69,4).Value = 34
$ws.Cells.Item(69,5).Value = 256
$ws.Cells.Item(69,8).Value = 10

# Row 70: Bulgaria -- rows 66-72 reshuffle -- Bulgaria jumps to row 70 with new stats
$ws.Cells.Item(70,1).Value = "Bulgaria"
$ws.Cells.Item(70,2).Value = 276
$ws.Cells.Item(70,3).Value = 12
$ws.Cells.Item(70,5).Value = 265
$ws.Cells.Item(70,6).Value = 8
$ws.Cells.Item(70,8).Value = 3

# Row 71: Marruecos -- rows 66-72 reshuffle -- shifts to row 71, stats unchanged
$ws.Cells.Item(71,1).Value = "Marruecos"
$ws.Cells.Item(71,2).Value = 275
$ws.Cells.Item(71,3).Value = 0
$ws.Cells.Item(71,4).Value = 8
$ws.Cells.Item(71,5).Value = 256
$ws.Cells.Item(71,6).Value = 1
$ws.Cells.Item(71,8).Value = 11

# Row 72: Taiwan -- rows 66-72 reshuffle -- shifts to row 72, stats unchanged
$ws.Cells.Item(72,1).Value = "Taiwan"
$ws.Cells.Item(72,2).Value = 267
$ws.Cells.Item(72,3).Value = 15
$ws.Cells.Item(72,4).Value = 30
$ws.Cells.Item(72,5).Value = 235
$ws.Cells.Item(72,6).Value = 0
$ws.Cells.Item(72,8).Value = 2

# Row 94: Kazajistan -- row 94/95 swap with Azerbaiyan -- Kazajistan overtakes it, new stats
$ws.Cells.Item(94,1).Value = "Kazajistan"
$ws.Cells.Item(94,2).Value = 124
$ws.Cells.Item(94,3).Value = 11
$ws.Cells.Item(94,4).Value = 2
$ws.Cells.Item(94,5).Value = 121
$ws.Cells.Item(94,6).Value = 0
$ws.Cells.Item(94,8).Value = 1

# Row 95: Azerbaiyan -- row 94/95 swap with Kazajistan -- Azerbaiyan drops a rank, stats unchanged
$ws.Cells.Item(95,1).Value = "Azerbaiyan"
$ws.Cells.Item(95,2).Value = 122
$ws.Cells.Item(95,3).Value = 0
$ws.Cells.Item(95,4).Value = 15
$ws.Cells.Item(95,5).Value = 104
$ws.Cells.Item(95,6).Value = 6
$ws.Cells.Item(95,8).Value = 3

# Row 107: Uzbekistan -- updated totals
$ws.Cells.Item(107,4).Value = 5
$ws.Cells.Item(107,5).Value = 77
$ws.Cells.Item(107,6).Value = 8
$ws.Cells.Item(107,7).Value = 1
$ws.Cells.Item(107,8).Value = 1

# Row 110: Georgia -- updated totals
$ws.Cells.Item(110,2).Value = 81
$ws.Cells.Item(110,3).Value = 2
$ws.Cells.Item(110,4).Value = 13

# Row 120: Consejo Danes para los Refugiados -- updated totals
$ws.Cells.Item(120,4).Value = 2
$ws.Cells.Item(120,5).Value = 46
